# Add a "Save" column (H) to the s_vals sheet, matching the header
# style used by the other header cells (B1:G1) and a value of 1 for H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the values first.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1

# Copy the formatting (bold, centered, bordered header style) from the
# existing "sum" header cell (G1) onto the new "Save" header cell (H1),
# without touching the value we just set.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
